$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F ("id"), shifting the old F.. (statuses_count etc.) to G.. onward.
# Column A (row labels: count/mean/std/min/25%/50%/75%/max) is untouched by this insert.
$ws.Columns("F:F").Insert()

# Re-apply header and data cell values for the refreshed descriptive-statistics table
$ws.Range("B1").Value2 = "default_profile"
$ws.Range("C1").Value2 = "favourites_count"
$ws.Range("D1").Value2 = "followers_count"
$ws.Range("E1").Value2 = "friends_count"
$ws.Range("F1").Value2 = "id"
$ws.Range("G1").Value2 = "statuses_count"
$ws.Range("H1").Value2 = "verified"
$ws.Range("I1").Value2 = "account_age_days"
$ws.Range("J1").Value2 = "ratio_statuses_count_per_age"
$ws.Range("K1").Value2 = "ratio_favorites_per_age"
$ws.Range("L1").Value2 = "ratio_friends_per_followers"
$ws.Range("M1").Value2 = "word_count"
$ws.Range("N1").Value2 = "char_count"
$ws.Range("O1").Value2 = "reputation"
$ws.Range("P1").Value2 = "description_word_count"
$ws.Range("Q1").Value2 = "description_character_count"
$ws.Range("R1").Value2 = "avg_word"
$ws.Range("B2").Value2 = 1200
$ws.Range("C2").Value2 = 1200
$ws.Range("D2").Value2 = 1200
$ws.Range("E2").Value2 = 1200
$ws.Range("F2").Value2 = 1200
$ws.Range("G2").Value2 = 1200
$ws.Range("H2").Value2 = 1200
$ws.Range("I2").Value2 = 1200
$ws.Range("J2").Value2 = 1200
$ws.Range("K2").Value2 = 1200
$ws.Range("L2").Value2 = 1200
$ws.Range("M2").Value2 = 1200
$ws.Range("N2").Value2 = 1200
$ws.Range("O2").Value2 = 1200
$ws.Range("P2").Value2 = 1200
$ws.Range("Q2").Value2 = 1200
$ws.Range("R2").Value2 = 1200
$ws.Range("B3").Value2 = 39.41666666666666
$ws.Range("C3").Value2 = 2.972588272546273
$ws.Range("D3").Value2 = 0.8132173887041124
$ws.Range("E3").Value2 = 0.2899707726986188
$ws.Range("F3").Value2 = 117659042009189000
$ws.Range("G3").Value2 = 0.8950400325647899
$ws.Range("H3").Value2 = 21.33333333333333
$ws.Range("I3").Value2 = 55.11941470378299
$ws.Range("J3").Value2 = 7.888849631600142
$ws.Range("K3").Value2 = 3.631691511453875
$ws.Range("L3").Value2 = 2.709531735934661
$ws.Range("M3").Value2 = 9.783333333333333
$ws.Range("N3").Value2 = 64.5325
$ws.Range("O3").Value2 = 63.31464341602851
$ws.Range("P3").Value2 = 9.484999999999999
$ws.Range("Q3").Value2 = 55.565
$ws.Range("R3").Value2 = 4.392518582122441
$ws.Range("B4").Value2 = 48.8874698617153
$ws.Range("C4").Value2 = 7.271878745159778
$ws.Range("D4").Value2 = 4.543025665767
$ws.Range("E4").Value2 = 3.512965866682032
$ws.Range("F4").Value2 = 292421216130287000
$ws.Range("G4").Value2 = 4.34661096472953
$ws.Range("H4").Value2 = 40.98319054374419
$ws.Range("I4").Value2 = 21.43569010288186
$ws.Range("J4").Value2 = 42.77300667365063
$ws.Range("K4").Value2 = 9.919048219479846
$ws.Range("L4").Value2 = 11.47763952068743
$ws.Range("M4").Value2 = 8.654891382502997
$ws.Range("N4").Value2 = 54.51399291254852
$ws.Range("O4").Value2 = 35.34234229180773
$ws.Range("P4").Value2 = 7.960061378382896
$ws.Range("Q4").Value2 = 46.97342717662045
$ws.Range("R4").Value2 = 6.463395584683599
$ws.Range("B5").Value2 = 0
$ws.Range("C5").Value2 = 0
$ws.Range("D5").Value2 = 0
$ws.Range("E5").Value2 = 0
$ws.Range("F5").Value2 = 418
$ws.Range("G5").Value2 = 0
$ws.Range("H5").Value2 = 0
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = 0
$ws.Range("K5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = 1
$ws.Range("N5").Value2 = 1
$ws.Range("O5").Value2 = 0
$ws.Range("P5").Value2 = 1
$ws.Range("Q5").Value2 = 1
$ws.Range("R5").Value2 = 0
$ws.Range("B6").Value2 = 0
$ws.Range("C6").Value2 = 0.09344266453942492
$ws.Range("D6").Value2 = 0.00007955839593494206
$ws.Range("E6").Value2 = 0.001968549369181749
$ws.Range("F6").Value2 = 94439953.25
$ws.Range("G6").Value2 = 0.0491718706595813
$ws.Range("H6").Value2 = 0
$ws.Range("I6").Value2 = 40.5085653104925
$ws.Range("J6").Value2 = 0.5013136288998358
$ws.Range("K6").Value2 = 0.1124325701065373
$ws.Range("L6").Value2 = 0.0007312739345898219
$ws.Range("M6").Value2 = 1.75
$ws.Range("N6").Value2 = 12
$ws.Range("O6").Value2 = 31.37138574907888
$ws.Range("P6").Value2 = 2
$ws.Range("Q6").Value2 = 10.75
$ws.Range("R6").Value2 = 2.438435538387252
$ws.Range("B7").Value2 = 0
$ws.Range("C7").Value2 = 0.5757364775920826
$ws.Range("D7").Value2 = 0.0008438011690069612
$ws.Range("E7").Value2 = 0.01842969495630499
$ws.Range("F7").Value2 = 377230939
$ws.Range("G7").Value2 = 0.1600520940434574
$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = 59.23982869379014
$ws.Range("J7").Value2 = 1.443625841750842
$ws.Range("K7").Value2 = 0.6632034452872573
$ws.Range("L7").Value2 = 0.4653561037105763
$ws.Range("M7").Value2 = 8
$ws.Range("N7").Value2 = 54
$ws.Range("O7").Value2 = 68.24302202821855
$ws.Range("P7").Value2 = 8
$ws.Range("Q7").Value2 = 47
$ws.Range("R7").Value2 = 3.846153846153847
$ws.Range("B8").Value2 = 100
$ws.Range("C8").Value2 = 2.347553155205024
$ws.Range("D8").Value2 = 0.03202165164869346
$ws.Range("E8").Value2 = 0.05518726334981938
$ws.Range("F8").Value2 = 2338591177.25
$ws.Range("G8").Value2 = 0.5945900119412246
$ws.Range("H8").Value2 = 0
$ws.Range("I8").Value2 = 73.29764453961455
$ws.Range("J8").Value2 = 5.358665313230251
$ws.Range("K8").Value2 = 2.585386326599655
$ws.Range("L8").Value2 = 2.187619047619048
$ws.Range("M8").Value2 = 16
$ws.Range("N8").Value2 = 113
$ws.Range("O8").Value2 = 99.92692604804158
$ws.Range("P8").Value2 = 16
$ws.Range("Q8").Value2 = 98
$ws.Range("R8").Value2 = 5.166586190246258
$ws.Range("B9").Value2 = 100
$ws.Range("C9").Value2 = 100
$ws.Range("D9").Value2 = 99.99999999999999
$ws.Range("E9").Value2 = 100
$ws.Range("F9").Value2 = 1118951206448181000
$ws.Range("G9").Value2 = 100
$ws.Range("H9").Value2 = 100
$ws.Range("I9").Value2 = 99.99999999999999
$ws.Range("J9").Value2 = 1062.034482758621
$ws.Range("K9").Value2 = 191.3500539374326
$ws.Range("L9").Value2 = 242.3333333333333
$ws.Range("M9").Value2 = 90
$ws.Range("N9").Value2 = 165
$ws.Range("O9").Value2 = 100
$ws.Range("P9").Value2 = 32
$ws.Range("Q9").Value2 = 156
$ws.Range("R9").Value2 = 100
